# Daily attendance processing - 2025-12-29 14:34:14
#
# For every row in the "Recorded By" column (G), the value is a
# comma-separated list of recorder identifiers (e.g. "System, someone@x.com").
# Normalize the list by swapping the first and last entries whenever the
# first entry is exactly "System" and the list doesn't already contain the
# admin@admin.com address (those rows are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    $parts = $val -split ", "

    if ($parts.Count -ge 2 -and $parts[0] -eq "System" -and -not $val.Contains("admin@admin.com")) {
        $lastIdx = $parts.Count - 1
        $tmp = $parts[0]
        $parts[0] = $parts[$lastIdx]
        $parts[$lastIdx] = $tmp
        $cell.Value = $parts -join ", "
    }
}
